# "Aggiornate curve di urgenza ed analisi sito esistente"
#
# Updates the "urgency curve" data points on Foglio1 (which drive the
# scatter chart) and renames the four chart series to their new labels.
# Also moves the active selection from B12 to B11 (and drops the stale
# topLeftCell scroll position that pointed at A2).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Foglio1")

# --- Updated urgency-curve coordinates -----------------------------------
# Row 6  (series "Prevenzione Errori", formerly "Accessibilità")
$ws.Range("A6").Value = 2
$ws.Range("B6").Value = 5

# Row 7  (series "Recognition~Recall", formerly "Convenzionalità")
$ws.Range("A7").Value = 4
$ws.Range("B7").Value = 3

# Row 8  (series "Flessibilità ed Efficienza", formerly "Credibilità")
$ws.Range("A8").Value = 2
$ws.Range("B8").Value = 4

# Row 9  (series "Estetica-Design Minimalista", formerly "Centralità dell'Utente")
$ws.Range("A9").Value = 5
$ws.Range("B9").Value = 2

# --- Rename the chart series ----------------------------------------------
$chartObj = $ws.ChartObjects().Item(1)
$chart = $chartObj.Chart
$series = $chart.SeriesCollection()

$series.Item(1).Name = "Prevenzione Errori"
$series.Item(2).Name = "Recognition~Recall"
$series.Item(3).Name = "Flessibilità ed Efficienza"
$series.Item(4).Name = "Estetica-Design Minimalista"

# --- Update the active selection ------------------------------------------
$null = $ws.Range("B11").Select()
